# Generate Report for Handback
# Updates the Overview status text, records the Handback datetime + file
# links on the zh-cn and de-de detail sheets, and widens a few columns so
# the new long file names / status text are readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# 1. Overview sheet: status flips from "Ready for handoff" to
#    "Handed back: in sync with en-US" for both locale columns/rows.
# ---------------------------------------------------------------------
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# Widen the two status columns so the longer text fits.
$overview.Range("E1").ColumnWidth = 29.9777050018311
$overview.Range("F1").ColumnWidth = 29.9777050018311

# ---------------------------------------------------------------------
# 2. zh-cn sheet: fill in Latest Target File / Latest Handback File /
#    Latest Handback DateTime for both rows, and hyperlink the new
#    "Latest Target File" cells the same way column A is linked.
# ---------------------------------------------------------------------
$zhcn.Range("I2").Value = "33ca47dd-6d15-476f-9e7b-e4b8e0d8eb0b.md"
$zhcn.Range("J2").Value = "33ca47dd-6d15-476f-9e7b-e4b8e0d8eb0b.79d56356936010aec9a1420b174b6b440aba3a64.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-31 12:55:26"

$zhcn.Range("I3").Value = "b8ada2d1-716b-400d-b079-415677aed1fe.md"
$zhcn.Range("J3").Value = "b8ada2d1-716b-400d-b079-415677aed1fe.294547b0747bfb98d966b055880843d121cbb00d.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-31 12:55:26"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1ba5fc72c72158f7ebe65e62c808b0754ae5ea0b/e2e/33ca47dd-6d15-476f-9e7b-e4b8e0d8eb0b.md", "", "", "33ca47dd-6d15-476f-9e7b-e4b8e0d8eb0b.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1ba5fc72c72158f7ebe65e62c808b0754ae5ea0b/e2e/b8ada2d1-716b-400d-b079-415677aed1fe.md", "", "", "b8ada2d1-716b-400d-b079-415677aed1fe.md")

$zhcn.Range("I2").Style = "HyperLink"
$zhcn.Range("I3").Style = "HyperLink"

$zhcn.Range("C1").ColumnWidth = 29.9777050018311
$zhcn.Range("I1").ColumnWidth = 40
$zhcn.Range("J1").ColumnWidth = 40

# ---------------------------------------------------------------------
# 3. de-de sheet: same shape of change as zh-cn, with de-de file names
#    and its own handback timestamp.
# ---------------------------------------------------------------------
$dede.Range("I2").Value = "33ca47dd-6d15-476f-9e7b-e4b8e0d8eb0b.md"
$dede.Range("J2").Value = "33ca47dd-6d15-476f-9e7b-e4b8e0d8eb0b.79d56356936010aec9a1420b174b6b440aba3a64.de-de.xlf"
$dede.Range("K2").Value = "2016-08-31 12:55:45"

$dede.Range("I3").Value = "b8ada2d1-716b-400d-b079-415677aed1fe.md"
$dede.Range("J3").Value = "b8ada2d1-716b-400d-b079-415677aed1fe.294547b0747bfb98d966b055880843d121cbb00d.de-de.xlf"
$dede.Range("K3").Value = "2016-08-31 12:55:45"

$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1ba5fc72c72158f7ebe65e62c808b0754ae5ea0b/e2e/33ca47dd-6d15-476f-9e7b-e4b8e0d8eb0b.md", "", "", "33ca47dd-6d15-476f-9e7b-e4b8e0d8eb0b.md")
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1ba5fc72c72158f7ebe65e62c808b0754ae5ea0b/e2e/b8ada2d1-716b-400d-b079-415677aed1fe.md", "", "", "b8ada2d1-716b-400d-b079-415677aed1fe.md")

$dede.Range("I2").Style = "HyperLink"
$dede.Range("I3").Style = "HyperLink"

$dede.Range("C1").ColumnWidth = 29.9777050018311
$dede.Range("I1").ColumnWidth = 40
$dede.Range("J1").ColumnWidth = 40
